$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.987.52"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.897.19"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8302"
$ws.Range("E5").Value = "  +4.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.89"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3273"
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.47"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07019"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08082"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7653"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "1.905.01"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.242"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.14"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "29.981.35"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.09"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.843"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.52"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007745"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "2.152.74"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.945"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1736"
$ws.Range("E25").Value = "  +24.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.254"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.32"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.091"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.357"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.512"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05860"
$ws.Range("E32").Value = "  +8.58%  "
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7306"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01916"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.776"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4439"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.43"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8566"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.855"
$ws.Range("E43").Value = "  -5.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.899"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.90"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.545"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.777"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "992.18"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "2.046.72"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("E51").Value = "  +0.52%  "
